$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right after the existing row 442 (i.e. at position 443),
# pushing the old row 443 down to row 447 and leaving row 442 in place
# (its values will be overwritten below with the newest week's data).
$ws.Rows.Item(443).Insert()
$ws.Rows.Item(443).Insert()
$ws.Rows.Item(443).Insert()
$ws.Rows.Item(443).Insert()

function Set-Row {
    param(
        [int]$r,
        [double]$a,
        [string]$b,
        [string]$c,
        [double]$d,
        [double]$e,
        [string]$f,
        [double]$g,
        [string]$h,
        [double]$i,
        [string]$j,
        [string]$k,
        [string]$l,
        [double]$m,
        [double]$n,
        [double]$o,
        [double]$p,
        [string]$q,
        [string]$rOrigen,
        [double]$s,
        [double]$t
    )
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rOrigen
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}

# Row 442: updated to the new week (Feria Lagunitas de Puerto Montt - Nectarín, August Red, Especial)
Set-Row 442 4 "Feria Lagunitas de Puerto Montt" "Los Lagos" 44628 10 "Fruta" 100103 "Frutos de hueso (carozo)" 100103006 "Nectarín" "August Red" "Especial" 250 21000 21000 21000 "$/caja 15 kilos empedrada" "Región de O'Higgins" 1400 15

# Row 443: new - August Red, Primera
Set-Row 443 4 "Feria Lagunitas de Puerto Montt" "Los Lagos" 44628 10 "Fruta" 100103 "Frutos de hueso (carozo)" 100103006 "Nectarín" "August Red" "Primera" 500 16000 17000 16500 "$/caja 15 kilos empedrada" "Región de O'Higgins" 1100 15

# Row 444: new - June Pearl, Especial
Set-Row 444 4 "Feria Lagunitas de Puerto Montt" "Los Lagos" 44628 10 "Fruta" 100103 "Frutos de hueso (carozo)" 100103006 "Nectarín" "June Pearl" "Especial" 250 21000 21000 21000 "$/caja 15 kilos empedrada" "Región de O'Higgins" 1400 15

# Row 445: new - June Pearl, Primera
Set-Row 445 4 "Feria Lagunitas de Puerto Montt" "Los Lagos" 44628 10 "Fruta" 100103 "Frutos de hueso (carozo)" 100103006 "Nectarín" "June Pearl" "Primera" 500 16000 17000 16500 "$/caja 15 kilos empedrada" "Región de O'Higgins" 1100 15

# Row 446: the previous week's Early John / Especial entry, re-written here because row 442
# (which originally held it) stayed put above the insertion point and was overwritten above.
Set-Row 446 4 "Feria Lagunitas de Puerto Montt" "Los Lagos" 44544 10 "Fruta" 100103 "Frutos de hueso (carozo)" 100103006 "Nectarín" "Early John" "Especial" 350 21000 21000 21000 "$/caja 15 kilos empedrada" "Región de O'Higgins" 1400 15

# Row 447 already contains the previous week's Early John / Primera data, which was shifted
# down automatically by the row insertions above, so no further changes are required there.
